# Refresh the cryptocurrency price/volume table (Coin, Link, Price, Volume(1h))
# with the latest scraped values from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All "Price" (column D) values are plain text in the source data (e.g. "26.836.29"
# uses dots as thousands separators, and values like "0.7170" rely on a trailing
# zero). Prefixing with a literal apostrophe forces Excel's quote-prefix text entry
# so the value round-trips as text instead of being auto-converted to a number.
$q = "'"

$ws.Range("D2").Value = $q + '26.851.00'
$ws.Range("E2").Value = '  -2.22%  '

$ws.Range("D3").Value = $q + '1.815.90'
$ws.Range("E3").Value = '  -1.22%  '

$ws.Range("D4").Value = $q + '1.008'
$ws.Range("E4").Value = '  -0.45%  '

$ws.Range("D5").Value = $q + '1.008'
$ws.Range("E5").Value = '  -0.29%  '

$ws.Range("D6").Value = $q + '308.50'
$ws.Range("E6").Value = '  -1.82%  '

$ws.Range("D7").Value = $q + '0.4609'
$ws.Range("E7").Value = '  -2.70%  '

$ws.Range("D8").Value = $q + '0.3641'
$ws.Range("E8").Value = '  -1.40%  '

$ws.Range("D9").Value = $q + '0.07217'
$ws.Range("E9").Value = '  -3.25%  '

$ws.Range("D10").Value = $q + '0.8571'
$ws.Range("E10").Value = '  -3.20%  '

$ws.Range("D11").Value = $q + '19.71'
$ws.Range("E11").Value = '  -3.60%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = $q + '0.07521'
$ws.Range("E12").Value = '  +2.50%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = $q + '1.811.18'
$ws.Range("E13").Value = '  -5.52%  '

$ws.Range("D14").Value = $q + '5.322'
$ws.Range("E14").Value = '  -2.37%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = $q + '6.517'
$ws.Range("E15").Value = '  -1.00%  '

$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = $q + '91.71'
$ws.Range("E16").Value = '  -1.70%  '

$ws.Range("E17").Value = '  +0.02%  '

$ws.Range("D18").Value = $q + '0.000008565'
$ws.Range("E18").Value = '  -2.86%  '

$ws.Range("E19").Value = '  -0.35%  '

$ws.Range("B20").Value = 'WrappedBTC'
$ws.Range("C20").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D20").Value = $q + '26.864.56'
$ws.Range("E20").Value = '  -2.26%  '

$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").Value = $q + '14.41'
$ws.Range("E21").Value = '  -2.65%  '

$ws.Range("D22").Value = $q + '5.140'
$ws.Range("E22").Value = '  -3.50%  '

$ws.Range("D23").Value = $q + '10.49'
$ws.Range("E23").Value = '  -1.91%  '

$ws.Range("D24").Value = $q + '2.028.31'
$ws.Range("E24").Value = '  -4.81%  '

$ws.Range("D25").Value = $q + '151.07'
$ws.Range("E25").Value = '  -0.72%  '

$ws.Range("E26").Value = '  -3.17%  '

$ws.Range("D27").Value = $q + '18.13'
$ws.Range("E27").Value = '  -2.72%  '

$ws.Range("D28").Value = $q + '2.061'
$ws.Range("E28").Value = '  -4.10%  '

$ws.Range("E29").Value = '  -3.06%  '

$ws.Range("D30").Value = $q + '115.12'
$ws.Range("E30").Value = '  -2.21%  '

$ws.Range("D31").Value = $q + '0.08853'

$ws.Range("D32").Value = $q + '2.956'
$ws.Range("E32").Value = '  +0.59%  '

$ws.Range("D33").Value = $q + '4.407'
$ws.Range("E33").Value = '  -3.31%  '

$ws.Range("D34").Value = $q + '1.131'
$ws.Range("E34").Value = '  -4.27%  '

$ws.Range("D35").Value = $q + '0.7170'
$ws.Range("E35").Value = '  -5.40%  '

$ws.Range("D36").Value = $q + '1.008'
$ws.Range("E36").Value = '  -0.44%  '

$ws.Range("E37").Value = '  -2.81%  '

$ws.Range("E38").Value = '  +1.43%  '

$ws.Range("D39").Value = $q + '0.05241'
$ws.Range("E39").Value = '  -1.68%  '

$ws.Range("D40").Value = $q + '0.01909'
$ws.Range("E40").Value = '  -2.40%  '

$ws.Range("D41").Value = $q + '2.920'
$ws.Range("E41").Value = '  -2.44%  '

$ws.Range("D42").Value = $q + '7.154'
$ws.Range("E42").Value = '  -2.36%  '

$ws.Range("D43").Value = $q + '0.5127'
$ws.Range("E43").Value = '  -3.87%  '

$ws.Range("D44").Value = $q + '0.1619'
$ws.Range("E44").Value = '  -2.50%  '

$ws.Range("D45").Value = $q + '8.179'

$ws.Range("D46").Value = $q + '0.4793'
$ws.Range("E46").Value = '  -2.43%  '

$ws.Range("D47").Value = $q + '1.008'
$ws.Range("E47").Value = '  -0.39%  '

$ws.Range("D48").Value = $q + '10.09'
$ws.Range("E48").Value = '  -4.03%  '

$ws.Range("D49").Value = $q + '102.84'
$ws.Range("E49").Value = '  -2.11%  '

$ws.Range("D50").Value = $q + '1.616'
$ws.Range("E50").Value = '  -3.73%  '

$ws.Range("D51").Value = $q + '0.06189'
$ws.Range("E51").Value = '  -1.97%  '
